# This workbook holds a long list of daily price records (one row per day)
# for "Pepino ensalada" at Femacal de La Calera. The edit inserts one new
# record row right before the old row 104, pushing the existing rows
# 104-188 down to 105-189 (dimension grows from A1:R188 to A1:R189).
#
# The new row's "constant" columns (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, Clasificacion) match every
# other row in the sheet; its variable columns (Fecha, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades) carry the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104; this shifts the former rows 104-188
# down to 105-189 and updates the sheet dimension automatically.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new record's values.
$ws.Cells.Item(104, 1).Value = 3
$ws.Cells.Item(104, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(104, 3).Value = "Coquimbo"
$ws.Cells.Item(104, 4).Value = 44447
$ws.Cells.Item(104, 5).Value = 5
$ws.Cells.Item(104, 6).Value = 100112043
$ws.Cells.Item(104, 7).Value = "Pepino ensalada"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 95
$ws.Cells.Item(104, 11).Value = 13000
$ws.Cells.Item(104, 12).Value = 14000
$ws.Cells.Item(104, 13).Value = 13474
$ws.Cells.Item(104, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 192
$ws.Cells.Item(104, 17).Value = 70
$ws.Cells.Item(104, 18).Value = "Hortaliza"
